$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''29.116.70'
$ws.Range('E2').Value = '''  +0.23%  '

$ws.Range('D3').Value = '''1.831.54'
$ws.Range('E3').Value = '''  -0.13%  '

$ws.Range('D4').Value = '''0.9993'
$ws.Range('E4').Value = '''  -0.04%  '

$ws.Range('D5').Value = '''243.18'
$ws.Range('E5').Value = '''  +0.09%  '

$ws.Range('D6').Value = '''0.6263'
$ws.Range('E6').Value = '''  -0.05%  '

$ws.Range('D7').Value = '''1.000'
$ws.Range('E7').Value = '''  -0.02%  '

$ws.Range('D8').Value = '''0.07510'
$ws.Range('E8').Value = '''  -0.90%  '

$ws.Range('D9').Value = '''0.2929'
$ws.Range('E9').Value = '''  +0.09%  '

$ws.Range('D10').Value = '''23.24'
$ws.Range('E10').Value = '''  +2.90%  '

$ws.Range('D11').Value = '''0.07704'
$ws.Range('E11').Value = '''  -0.46%  '

$ws.Range('D12').Value = '''1.829.38'
$ws.Range('E12').Value = '''  -0.42%  '

$ws.Range('D13').Value = '''5.028'
$ws.Range('E13').Value = '''  +1.27%  '

$ws.Range('D14').Value = '''0.6692'
$ws.Range('E14').Value = '''  +0.67%  '

$ws.Range('D15').Value = '''82.82'
$ws.Range('E15').Value = '''  -0.23%  '

$ws.Range('D16').Value = '''0.000009370'
$ws.Range('E16').Value = '''  -7.01%  '

$ws.Range('D17').Value = '''5.996'
$ws.Range('E17').Value = '''  -1.22%  '

$ws.Range('D18').Value = '''29.101.14'
$ws.Range('E18').Value = '''  +0.10%  '

$ws.Range('D19').Value = '''2.075.65'
$ws.Range('E19').Value = '''  -0.48%  '

$ws.Range('D20').Value = '''12.61'
$ws.Range('E20').Value = '''  +1.69%  '

$ws.Range('D21').Value = '''222.89'
$ws.Range('E21').Value = '''  -1.72%  '

$ws.Range('D22').Value = '''1.002'

$ws.Range('D23').Value = '''7.157'
$ws.Range('E23').Value = '''  -0.76%  '

$ws.Range('D24').Value = '''1.000'
$ws.Range('E24').Value = '''  +0.00%  '

$ws.Range('D25').Value = '''160.32'
$ws.Range('E25').Value = '''  +0.46%  '

$ws.Range('D26').Value = '''0.1397'
$ws.Range('E26').Value = '''  +0.82%  '

$ws.Range('D27').Value = '''8.506'
$ws.Range('E27').Value = '''  -0.03%  '

$ws.Range('D28').Value = '''17.89'

$ws.Range('D29').Value = '''1.491'
$ws.Range('E29').Value = '''  -0.14%  '

$ws.Range('D30').Value = '''0.05822'
$ws.Range('E30').Value = '''  +10.88%  '

$ws.Range('D31').Value = '''4.161'
$ws.Range('E31').Value = '''  +1.54%  '

$ws.Range('D32').Value = '''4.126'
$ws.Range('E32').Value = '''  +2.89%  '

$ws.Range('D33').Value = '''1.209'
$ws.Range('E33').Value = '''  +1.28%  '

$ws.Range('D34').Value = '''0.7427'
$ws.Range('E34').Value = '''  +1.04%  '

$ws.Range('D35').Value = '''1.830'
$ws.Range('E35').Value = '''  -0.64%  '

$ws.Range('D36').Value = '''1.140'
$ws.Range('E36').Value = '''  +0.31%  '

$ws.Range('D37').Value = '''2.667'
$ws.Range('E37').Value = '''  -0.84%  '

$ws.Range('D38').Value = '''1.231.78'
$ws.Range('E38').Value = '''  -0.93%  '

$ws.Range('D39').Value = '''2.766'
$ws.Range('E39').Value = '''  +0.08%  '

$ws.Range('E40').Value = '''  -0.30%  '

$ws.Range('D41').Value = '''6.491'
$ws.Range('E41').Value = '''  +1.79%  '

$ws.Range('D42').Value = '''0.8927'
$ws.Range('E42').Value = '''  -0.65%  '

$ws.Range('D43').Value = '''1.000'
$ws.Range('E43').Value = '''  +0.01%  '

$ws.Range('D44').Value = '''102.27'
$ws.Range('E44').Value = '''  +0.06%  '

$ws.Range('B45').Value = 'BabyDogeCoin'
$ws.Range('C45').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D45').Value = '''0.00000000126'
$ws.Range('E45').Value = '''  +1.66%  '

$ws.Range('B46').Value = 'RocketPoolETH'
$ws.Range('C46').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D46').Value = '''1.976.44'
$ws.Range('E46').Value = '''  -0.45%  '

$ws.Range('D47').Value = '''66.01'
$ws.Range('E47').Value = '''  +2.82%  '

$ws.Range('D48').Value = '''0.5088'
$ws.Range('E48').Value = '''  -0.50%  '

$ws.Range('B49').Value = 'TheSandbox'
$ws.Range('C49').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D49').Value = '''0.4069'
$ws.Range('E49').Value = '''  +0.76%  '

$ws.Range('B50').Value = 'XinFinNetwork'
$ws.Range('C50').Value = 'https://coinranking.com/coin/77jGXSqWJ1ofG+xinfinnetwork-xdc'
$ws.Range('D50').Value = '''0.07467'
$ws.Range('E50').Value = '''  +12.13%  '

$ws.Range('D51').Value = '''9.011'
$ws.Range('E51').Value = '''  +1.44%  '
